$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix start week dates issue: milestone labels were showing the wrong text ---
$ws.Range("C5").Value = "M1"
$ws.Range("C8").Value = "M2"

# --- Add a thin border around each task's summary row (name/activity + start/end dates) ---
$borderedCells = @("B5", "C5", "D5", "E5", "B8", "C8", "D8", "E8")
foreach ($addr in $borderedCells) {
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).Borders.Weight = 2
}

# --- Highlight the overall task span (summary row week cells) in green so it is
#     visually distinct from the per-subtask week highlight, which stays orange ---
$greenCells = @("F5", "G5", "H8", "I8")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 5417010
}
